$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID19_TIMESERIESDATA")

# --- New block appended for the 07-10-2020 report (rows 721-756) ---
# Row 721 repeats the column headers, just like every earlier daily block
# (e.g. rows 1, 433, 469, ... 685) that starts a new day of data.
$ws.Cells.Item(721,1).Value = "States/UT"
$ws.Cells.Item(721,2).Value = "Active Cases"
$ws.Cells.Item(721,3).Value = "Active Cases Since Yesterday"
$ws.Cells.Item(721,4).Value = "Recovered Cases"
$ws.Cells.Item(721,5).Value = "Recovered Cases Since Yesterday"
$ws.Cells.Item(721,6).Value = "Deceased Cases"
$ws.Cells.Item(721,7).Value = "Deceased Cases Since Yesterday"
$ws.Cells.Item(721,8).Value = "Date"

# Give the new header row the same bold / boxed / centered look as the
# most recent previous header row (row 433) instead of leaving it unformatted.
$ws.Range("A433:H433").Copy()
$ws.Range("A721:H721").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: one per State/UT, in the same column order as the rest of the sheet.
# Column H holds the report date as literal text (e.g. "07-10-2020"), matching the
# existing cells, so force a Text format first to stop Excel auto-converting it to a date serial.
$ws.Range("H722:H756").NumberFormat = "@"

$ws.Cells.Item(722,1).Value = "Andaman and Nicobar Islands"
$ws.Cells.Item(722,2).Value = 180
$ws.Cells.Item(722,3).Value = -6
$ws.Cells.Item(722,4).Value = 3678
$ws.Cells.Item(722,5).Value = 19
$ws.Cells.Item(722,6).Value = 54
$ws.Cells.Item(722,7).Value = 0
$ws.Cells.Item(722,8).Value = "07-10-2020"
$ws.Cells.Item(723,1).Value = "Andhra Pradesh"
$ws.Cells.Item(723,2).Value = 50776
$ws.Cells.Item(723,3).Value = -284
$ws.Cells.Item(723,4).Value = 672479
$ws.Cells.Item(723,5).Value = 6046
$ws.Cells.Item(723,6).Value = 6052
$ws.Cells.Item(723,7).Value = 33
$ws.Cells.Item(723,8).Value = "07-10-2020"
$ws.Cells.Item(724,1).Value = "Arunachal Pradesh"
$ws.Cells.Item(724,2).Value = 3022
$ws.Cells.Item(724,3).Value = 33
$ws.Cells.Item(724,4).Value = 7965
$ws.Cells.Item(724,5).Value = 190
$ws.Cells.Item(724,6).Value = 20
$ws.Cells.Item(724,7).Value = 1
$ws.Cells.Item(724,8).Value = "07-10-2020"
$ws.Cells.Item(725,1).Value = "Assam"
$ws.Cells.Item(725,2).Value = 33047
$ws.Cells.Item(725,3).Value = -420
$ws.Cells.Item(725,4).Value = 155077
$ws.Cells.Item(725,5).Value = 1586
$ws.Cells.Item(725,6).Value = 778
$ws.Cells.Item(725,7).Value = 18
$ws.Cells.Item(725,8).Value = "07-10-2020"
$ws.Cells.Item(726,1).Value = "Bihar"
$ws.Cells.Item(726,2).Value = 11420
$ws.Cells.Item(726,3).Value = -103
$ws.Cells.Item(726,4).Value = 178395
$ws.Cells.Item(726,5).Value = 1400
$ws.Cells.Item(726,6).Value = 925
$ws.Cells.Item(726,7).Value = 1
$ws.Cells.Item(726,8).Value = "07-10-2020"
$ws.Cells.Item(727,1).Value = "Chandigarh"
$ws.Cells.Item(727,2).Value = 1492
$ws.Cells.Item(727,3).Value = -112
$ws.Cells.Item(727,4).Value = 11035
$ws.Cells.Item(727,5).Value = 238
$ws.Cells.Item(727,6).Value = 180
$ws.Cells.Item(727,7).Value = 3
$ws.Cells.Item(727,8).Value = "07-10-2020"
$ws.Cells.Item(728,1).Value = "Chhattisgarh"
$ws.Cells.Item(728,2).Value = 27238
$ws.Cells.Item(728,3).Value = -619
$ws.Cells.Item(728,4).Value = 100551
$ws.Cells.Item(728,5).Value = 3484
$ws.Cells.Item(728,6).Value = 1104
$ws.Cells.Item(728,7).Value = 23
$ws.Cells.Item(728,8).Value = "07-10-2020"
$ws.Cells.Item(729,1).Value = "Dadra and Nagar Haveli and Daman and Diu"
$ws.Cells.Item(729,2).Value = 101
$ws.Cells.Item(729,3).Value = 2
$ws.Cells.Item(729,4).Value = 3000
$ws.Cells.Item(729,5).Value = 9
$ws.Cells.Item(729,6).Value = 2
$ws.Cells.Item(729,7).Value = 0
$ws.Cells.Item(729,8).Value = "07-10-2020"
$ws.Cells.Item(730,1).Value = "Delhi"
$ws.Cells.Item(730,2).Value = 22720
$ws.Cells.Item(730,3).Value = -360
$ws.Cells.Item(730,4).Value = 266935
$ws.Cells.Item(730,5).Value = 2997
$ws.Cells.Item(730,6).Value = 5581
$ws.Cells.Item(730,7).Value = 39
$ws.Cells.Item(730,8).Value = "07-10-2020"
$ws.Cells.Item(731,1).Value = "Goa"
$ws.Cells.Item(731,2).Value = 4720
$ws.Cells.Item(731,3).Value = -83
$ws.Cells.Item(731,4).Value = 31050
$ws.Cells.Item(731,5).Value = 594
$ws.Cells.Item(731,6).Value = 468
$ws.Cells.Item(731,7).Value = 8
$ws.Cells.Item(731,8).Value = "07-10-2020"
$ws.Cells.Item(732,1).Value = "Gujarat"
$ws.Cells.Item(732,2).Value = 16570
$ws.Cells.Item(732,3).Value = -148
$ws.Cells.Item(732,4).Value = 125111
$ws.Cells.Item(732,5).Value = 1473
$ws.Cells.Item(732,6).Value = 3519
$ws.Cells.Item(732,7).Value = 10
$ws.Cells.Item(732,8).Value = "07-10-2020"
$ws.Cells.Item(733,1).Value = "Haryana"
$ws.Cells.Item(733,2).Value = 11320
$ws.Cells.Item(733,3).Value = -502
$ws.Cells.Item(733,4).Value = 123286
$ws.Cells.Item(733,5).Value = 1690
$ws.Cells.Item(733,6).Value = 1509
$ws.Cells.Item(733,7).Value = 18
$ws.Cells.Item(733,8).Value = "07-10-2020"
$ws.Cells.Item(734,1).Value = "Himachal Pradesh"
$ws.Cells.Item(734,2).Value = 3136
$ws.Cells.Item(734,3).Value = -20
$ws.Cells.Item(734,4).Value = 12918
$ws.Cells.Item(734,5).Value = 265
$ws.Cells.Item(734,6).Value = 229
$ws.Cells.Item(734,7).Value = 5
$ws.Cells.Item(734,8).Value = "07-10-2020"
$ws.Cells.Item(735,1).Value = "Jammu and Kashmir"
$ws.Cells.Item(735,2).Value = 13712
$ws.Cells.Item(735,3).Value = -984
$ws.Cells.Item(735,4).Value = 65496
$ws.Cells.Item(735,5).Value = 1706
$ws.Cells.Item(735,6).Value = 1268
$ws.Cells.Item(735,7).Value = 16
$ws.Cells.Item(735,8).Value = "07-10-2020"
$ws.Cells.Item(736,1).Value = "Jharkhand"
$ws.Cells.Item(736,2).Value = 10027
$ws.Cells.Item(736,3).Value = -409
$ws.Cells.Item(736,4).Value = 78089
$ws.Cells.Item(736,5).Value = 1246
$ws.Cells.Item(736,6).Value = 757
$ws.Cells.Item(736,7).Value = 10
$ws.Cells.Item(736,8).Value = "07-10-2020"
$ws.Cells.Item(737,1).Value = "Karnataka"
$ws.Cells.Item(737,2).Value = 115170
$ws.Cells.Item(737,3).Value = -326
$ws.Cells.Item(737,4).Value = 533074
$ws.Cells.Item(737,5).Value = 10228
$ws.Cells.Item(737,6).Value = 9461
$ws.Cells.Item(737,7).Value = 91
$ws.Cells.Item(737,8).Value = "07-10-2020"
$ws.Cells.Item(738,1).Value = "Kerala"
$ws.Cells.Item(738,2).Value = 87823
$ws.Cells.Item(738,3).Value = 2865
$ws.Cells.Item(738,4).Value = 154092
$ws.Cells.Item(738,5).Value = 4981
$ws.Cells.Item(738,6).Value = 884
$ws.Cells.Item(738,7).Value = 25
$ws.Cells.Item(738,8).Value = "07-10-2020"
$ws.Cells.Item(739,1).Value = "Ladakh"
$ws.Cells.Item(739,2).Value = 1195
$ws.Cells.Item(739,3).Value = 29
$ws.Cells.Item(739,4).Value = 3464
$ws.Cells.Item(739,5).Value = 50
$ws.Cells.Item(739,6).Value = 61
$ws.Cells.Item(739,7).Value = 0
$ws.Cells.Item(739,8).Value = "07-10-2020"
$ws.Cells.Item(740,1).Value = "Madhya Pradesh"
$ws.Cells.Item(740,2).Value = 18141
$ws.Cells.Item(740,3).Value = -616
$ws.Cells.Item(740,4).Value = 118039
$ws.Cells.Item(740,5).Value = 2161
$ws.Cells.Item(740,6).Value = 2488
$ws.Cells.Item(740,7).Value = 25
$ws.Cells.Item(740,8).Value = "07-10-2020"
$ws.Cells.Item(741,1).Value = "Maharashtra"
$ws.Cells.Item(741,2).Value = 247468
$ws.Cells.Item(741,3).Value = -5253
$ws.Cells.Item(741,4).Value = 1179726
$ws.Cells.Item(741,5).Value = 17141
$ws.Cells.Item(741,6).Value = 38717
$ws.Cells.Item(741,7).Value = 370
$ws.Cells.Item(741,8).Value = "07-10-2020"
$ws.Cells.Item(742,1).Value = "Manipur"
$ws.Cells.Item(742,2).Value = 2680
$ws.Cells.Item(742,3).Value = -16
$ws.Cells.Item(742,4).Value = 9482
$ws.Cells.Item(742,5).Value = 148
$ws.Cells.Item(742,6).Value = 78
$ws.Cells.Item(742,7).Value = 3
$ws.Cells.Item(742,8).Value = "07-10-2020"
$ws.Cells.Item(743,1).Value = "Meghalaya"
$ws.Cells.Item(743,2).Value = 2371
$ws.Cells.Item(743,3).Value = 154
$ws.Cells.Item(743,4).Value = 4606
$ws.Cells.Item(743,5).Value = 115
$ws.Cells.Item(743,6).Value = 60
$ws.Cells.Item(743,7).Value = 1
$ws.Cells.Item(743,8).Value = "07-10-2020"
$ws.Cells.Item(744,1).Value = "Mizoram"
$ws.Cells.Item(744,2).Value = 261
$ws.Cells.Item(744,3).Value = -30
$ws.Cells.Item(744,4).Value = 1887
$ws.Cells.Item(744,5).Value = 50
$ws.Cells.Item(744,6).Value = 0
$ws.Cells.Item(744,7).Value = 0
$ws.Cells.Item(744,8).Value = "07-10-2020"
$ws.Cells.Item(745,1).Value = "Nagaland"
$ws.Cells.Item(745,2).Value = 1185
$ws.Cells.Item(745,3).Value = 30
$ws.Cells.Item(745,4).Value = 5460
$ws.Cells.Item(745,5).Value = 38
$ws.Cells.Item(745,6).Value = 17
$ws.Cells.Item(745,7).Value = 0
$ws.Cells.Item(745,8).Value = "07-10-2020"
$ws.Cells.Item(746,1).Value = "Odisha"
$ws.Cells.Item(746,2).Value = 26846
$ws.Cells.Item(746,3).Value = -1160
$ws.Cells.Item(746,4).Value = 210217
$ws.Cells.Item(746,5).Value = 3817
$ws.Cells.Item(746,6).Value = 940
$ws.Cells.Item(746,7).Value = 16
$ws.Cells.Item(746,8).Value = "07-10-2020"
$ws.Cells.Item(747,1).Value = "Puducherry"
$ws.Cells.Item(747,2).Value = 4522
$ws.Cells.Item(747,3).Value = 9
$ws.Cells.Item(747,4).Value = 24614
$ws.Cells.Item(747,5).Value = 393
$ws.Cells.Item(747,6).Value = 546
$ws.Cells.Item(747,7).Value = 3
$ws.Cells.Item(747,8).Value = "07-10-2020"
$ws.Cells.Item(748,1).Value = "Punjab"
$ws.Cells.Item(748,2).Value = 11982
$ws.Cells.Item(748,3).Value = -913
$ws.Cells.Item(748,4).Value = 104355
$ws.Cells.Item(748,5).Value = 1707
$ws.Cells.Item(748,6).Value = 3679
$ws.Cells.Item(748,7).Value = 38
$ws.Cells.Item(748,8).Value = "07-10-2020"
$ws.Cells.Item(749,1).Value = "Rajasthan"
$ws.Cells.Item(749,2).Value = 21294
$ws.Cells.Item(749,3).Value = 79
$ws.Cells.Item(749,4).Value = 125448
$ws.Cells.Item(749,5).Value = 2027
$ws.Cells.Item(749,6).Value = 1574
$ws.Cells.Item(749,7).Value = 15
$ws.Cells.Item(749,8).Value = "07-10-2020"
$ws.Cells.Item(750,1).Value = "Sikkim"
$ws.Cells.Item(750,2).Value = 580
$ws.Cells.Item(750,3).Value = -18
$ws.Cells.Item(750,4).Value = 2587
$ws.Cells.Item(750,5).Value = 40
$ws.Cells.Item(750,6).Value = 49
$ws.Cells.Item(750,7).Value = 3
$ws.Cells.Item(750,8).Value = "07-10-2020"
$ws.Cells.Item(751,1).Value = "Tamil Nadu"
$ws.Cells.Item(751,2).Value = 45279
$ws.Cells.Item(751,3).Value = -602
$ws.Cells.Item(751,4).Value = 575212
$ws.Cells.Item(751,5).Value = 5548
$ws.Cells.Item(751,6).Value = 9917
$ws.Cells.Item(751,7).Value = 71
$ws.Cells.Item(751,8).Value = "07-10-2020"
$ws.Cells.Item(752,1).Value = "Telengana"
$ws.Cells.Item(752,2).Value = 26551
$ws.Cells.Item(752,3).Value = -93
$ws.Cells.Item(752,4).Value = 177008
$ws.Cells.Item(752,5).Value = 2239
$ws.Cells.Item(752,6).Value = 1189
$ws.Cells.Item(752,7).Value = 8
$ws.Cells.Item(752,8).Value = "07-10-2020"
$ws.Cells.Item(753,1).Value = "Tripura"
$ws.Cells.Item(753,2).Value = 4621
$ws.Cells.Item(753,3).Value = -255
$ws.Cells.Item(753,4).Value = 22623
$ws.Cells.Item(753,5).Value = 492
$ws.Cells.Item(753,6).Value = 301
$ws.Cells.Item(753,7).Value = 0
$ws.Cells.Item(753,8).Value = "07-10-2020"
$ws.Cells.Item(754,1).Value = "Uttarakhand"
$ws.Cells.Item(754,2).Value = 8414
$ws.Cells.Item(754,3).Value = -287
$ws.Cells.Item(754,4).Value = 43238
$ws.Cells.Item(754,5).Value = 617
$ws.Cells.Item(754,6).Value = 677
$ws.Cells.Item(754,7).Value = 8
$ws.Cells.Item(754,8).Value = "07-10-2020"
$ws.Cells.Item(755,1).Value = "Uttar Pradesh"
$ws.Cells.Item(755,2).Value = 44031
$ws.Cells.Item(755,3).Value = -993
$ws.Cells.Item(755,4).Value = 370753
$ws.Cells.Item(755,5).Value = 4432
$ws.Cells.Item(755,6).Value = 6153
$ws.Cells.Item(755,7).Value = 61
$ws.Cells.Item(755,8).Value = "07-10-2020"
$ws.Cells.Item(756,1).Value = "West Bengal"
$ws.Cells.Item(756,2).Value = 27988
$ws.Cells.Item(756,3).Value = 271
$ws.Cells.Item(756,4).Value = 243743
$ws.Cells.Item(756,5).Value = 3036
$ws.Cells.Item(756,6).Value = 5318
$ws.Cells.Item(756,7).Value = 63
$ws.Cells.Item(756,8).Value = "07-10-2020"

# The source cells for dates elsewhere in the sheet carry no explicit number
# format (plain inline text) - drop the Text format now that the literal
# string is safely stored, so the new cells match the rest of the column.
$ws.Range("H722:H756").ClearFormats()
